$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Remove existing hyperlinks in the data area so we can rebuild cleanly
$ws.Range("A1:H100").Hyperlinks.Delete()

# Clear old data rows (2-8) entirely (contents + formatting)
$ws.Range("A2:H20").Clear()

# Row 2
$ws.Range("A2").Value2 = "2026-01-09 18:29:08"
$ws.Range("B2").Value2 = "製造業向け図面自動生成システムの開発・ツール化を支援してくださるエンジニア募集(AI/バックエンド)"
$ws.Range("C2").Value2 = "システム開発"
$ws.Range("D2").Value2 = "300,000 円 ~ 500,000 円 / 固定"
$ws.Range("E2").Value2 = "期限情報なし"
$ws.Range("G2").Value2 = 435
$ws.Range("H2").Value2 = "🔥AI,Ai ◆ツール,開発"
$ws.Range("F2").Value2 = "https://www.lancers.jp/work/detail/5460562"
$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.lancers.jp/work/detail/5460562")

# Row 3
$ws.Range("A3").Value2 = "2026-01-09 18:29:08"
$ws.Range("B3").Value2 = "【AI×自動化】管理画面のデータ監視・チャットログのAI解析・LINE通知システムの構築依頼"
$ws.Range("C3").Value2 = "システム開発"
$ws.Range("D3").Value2 = "100,000 円 ~ 200,000 円 / 固定"
$ws.Range("E3").Value2 = "期限情報なし"
$ws.Range("G3").Value2 = 413
$ws.Range("H3").Value2 = "🔥AI,Ai ◆自動化 ◇管理"
$ws.Range("F3").Value2 = "https://www.lancers.jp/work/detail/5468493"
$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.lancers.jp/work/detail/5468493")

# Row 4
$ws.Range("A4").Value2 = "2026-01-09 18:29:08"
$ws.Range("B4").Value2 = "施設管理・現場業務向け チェックリスト業務の自動化・報告書作成システム開発エンジニア募集"
$ws.Range("C4").Value2 = "システム開発"
$ws.Range("D4").Value2 = "300,000 円 ~ 500,000 円 / 固定"
$ws.Range("E4").Value2 = "期限情報なし"
$ws.Range("G4").Value2 = 220
$ws.Range("H4").Value2 = "◆開発,システム開発 ◇管理"
$ws.Range("F4").Value2 = "https://www.lancers.jp/work/detail/5460563"
$ws.Hyperlinks.Add($ws.Range("F4"), "https://www.lancers.jp/work/detail/5460563")

# Row 5
$ws.Range("A5").Value2 = "2026-01-09 18:29:08"
$ws.Range("B5").Value2 = "初回 急募 自動カートインツール 開発のプロフェッショナルを探しています"
$ws.Range("C5").Value2 = "システム開発"
$ws.Range("D5").Value2 = "10,000 円 ~ 20,000 円 / 固定"
$ws.Range("E5").Value2 = "期限情報なし"
$ws.Range("G5").Value2 = 120
$ws.Range("H5").Value2 = "◆ツール,開発"
$ws.Range("F5").Value2 = "https://www.lancers.jp/work/detail/5467745"
$ws.Hyperlinks.Add($ws.Range("F5"), "https://www.lancers.jp/work/detail/5467745")

# Row 6
$ws.Range("A6").Value2 = "2026-01-09 18:29:08"
$ws.Range("B6").Value2 = "初回 スマホとの距離を見直す、シンプルなiOSアプリ開発"
$ws.Range("C6").Value2 = "システム開発"
$ws.Range("D6").Value2 = "300,000 円 ~ 500,000 円 / 固定"
$ws.Range("E6").Value2 = "期限情報なし"
$ws.Range("G6").Value2 = 100
$ws.Range("H6").Value2 = "◆開発 ◇アプリ"
$ws.Range("F6").Value2 = "https://www.lancers.jp/work/detail/5468441"
$ws.Hyperlinks.Add($ws.Range("F6"), "https://www.lancers.jp/work/detail/5468441")

# Row 7
$ws.Range("A7").Value2 = "2026-01-09 18:29:08"
$ws.Range("B7").Value2 = "【緊急募集】動画解析アプリ開発のプロフェッショナル"
$ws.Range("C7").Value2 = "システム開発"
$ws.Range("D7").Value2 = "20,000 円 ~ 50,000 円 / 固定"
$ws.Range("E7").Value2 = "期限情報なし"
$ws.Range("G7").Value2 = 88
$ws.Range("H7").Value2 = "◆開発 ◇アプリ"
$ws.Range("F7").Value2 = "https://www.lancers.jp/work/detail/5467910"
$ws.Hyperlinks.Add($ws.Range("F7"), "https://www.lancers.jp/work/detail/5467910")

# Row 8
$ws.Range("A8").Value2 = "2026-01-09 18:29:08"
$ws.Range("B8").Value2 = "【急募】BtoB向け越境ECプラットフォーム開発のパートナー募集"
$ws.Range("C8").Value2 = "システム開発"
$ws.Range("D8").Value2 = "3,000,000 円 ~ 5,000,000 円 / 固定"
$ws.Range("E8").Value2 = "期限情報なし"
$ws.Range("G8").Value2 = 75
$ws.Range("H8").Value2 = "◆開発"
$ws.Range("F8").Value2 = "https://www.lancers.jp/work/detail/5468347"
$ws.Hyperlinks.Add($ws.Range("F8"), "https://www.lancers.jp/work/detail/5468347")

# Row 9
$ws.Range("A9").Value2 = "2026-01-09 18:29:08"
$ws.Range("B9").Value2 = "【急募】大手保険システム会社でのPJ推進支援(PM・PL経験者募集/都内常駐)"
$ws.Range("C9").Value2 = "システム開発"
$ws.Range("D9").Value2 = "500,000 円 ~ 1,000,000 円 / 固定"
$ws.Range("E9").Value2 = "期限情報なし"
$ws.Range("G9").Value2 = 40
$ws.Range("F9").Value2 = "https://www.lancers.jp/work/detail/5467981"
$ws.Hyperlinks.Add($ws.Range("F9"), "https://www.lancers.jp/work/detail/5467981")

# Row 10
$ws.Range("A10").Value2 = "2026-01-09 18:29:08"
$ws.Range("B10").Value2 = "限定公開 限定公開の仕事"
$ws.Range("C10").Value2 = "システム開発"
$ws.Range("D10").Value2 = "5,000,000 円 ~ / 固定"
$ws.Range("E10").Value2 = "期限情報なし"
$ws.Range("G10").Value2 = 25
$ws.Range("F10").Value2 = "https://www.lancers.jp/work/detail/5467882"
$ws.Hyperlinks.Add($ws.Range("F10"), "https://www.lancers.jp/work/detail/5467882")

# Adjust column H width (16 -> 17) per target layout
$ws.Columns.Item(8).ColumnWidth = 16.14
